$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the row for "Control " (row 9) entirely - it was removed from the tracker.
$ws.Rows.Item(9).Delete()

# 2. Update Status (column B) for several games from "Completo" to their new values
#    (row numbers below are the POST-delete row numbers).
$ws.Range("B2").Value = "Platinado"    # Astro's Bot
$ws.Range("B5").Value = "Zerado"       # Bendy and the Ink Machine
$ws.Range("B7").Value = "Zerado"       # Celeste (PC - Steam / Platinar row)
$ws.Range("B11").Value = "Zerado"      # Days Gone
$ws.Range("B13").Value = "Zerado"      # Elden Ring
$ws.Range("B73").Value = "Zerado"      # Hollow Knight (PC - Xbox row)
$ws.Range("B74").Value = "Zerado"      # GTA Vice City

# 3. Highlight the "Crash Bandicoot" row (now row 9) with an underline font style,
#    including two extra (empty) styled cells to the right (E9:F9).
$rng = $ws.Range("A9:F9")
$rng.Font.Underline = $true

# 4. Restore view state: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("B76").Select()
